$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "prefix|Emissions|BC|Harmonized"
$ws.Range("D3").Value = "prefix|Emissions|BC|sector1|Harmonized"
$ws.Range("D4").Value = "prefix|Emissions|BC|sector2|Harmonized"
$ws.Range("D5").Value = "prefix|Emissions|BC|Harmonized"
$ws.Range("D6").Value = "prefix|Emissions|BC|sector1|Harmonized"
$ws.Range("D7").Value = "prefix|Emissions|BC|sector2|Harmonized"
